$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.990.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.19%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.423.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.36%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'572.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'163.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.422.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.49%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.553"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.52%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.42%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.49%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.424"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.43%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.013.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.61%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.23%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'27.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.81%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -4.86%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'64.045.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.16%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.365.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.10%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.11%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.46%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'378.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.82%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.13%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'71.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.28%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.517"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.52%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.32%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.40%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.37%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.15%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.58%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.06%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'22.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.16%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'7.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.71%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -3.49%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.72%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.857"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +10.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -3.97%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0730"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.799.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.50%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'25.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'42.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.07%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'26.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.23%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -2.71%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.82%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +8.77%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'329.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.37%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -3.40%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'6.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.77%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -2.41%  "
$ws.Range("E51").Style = "Normal"
